$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the two time-slot values in column C
$ws.Range("C2").Value = "9:05-9:10"
$ws.Range("C3").Value = "9:10-9:15"

# Update the active cell / selection shown in the sheet view
$ws.Range("C11").Select()
